$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Property1" to "DataNode" (unify DataNode/DataTable/Entity concepts)
$ws.Name = "DataNode"

# Update the active selection on the sheet to C38
$ws.Range("C38").Select()
